$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine last used row based on the "statut" column (A)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Map old emoji -> new emoji (statut column, A)
$emojiMap = @{
    "🟥" = "📕"
    "⬛" = "📘"
    "🟩" = "📗"
    "🟧" = "📙"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $val = $cellA.Value()
    if ($emojiMap.ContainsKey($val)) {
        $cellA.Value = $emojiMap[$val]
    }

    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value() -eq "noir") {
        $cellB.Value = "bleu"
    }
}
